# Update "想去人数" (F column) values for specific rows on both the
# "展览" and "全部类型" worksheets, which contain duplicate data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 200
    3  = 242
    7  = 6005
    15 = 387
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
